# Generate Report for Handback
# Updates the "zh-cn" and "de-de" localization-status sheets with the
# newly generated handback info for the ed324dd0-... file (row 8):
#   - adds a hyperlinked handback file name in column I
#   - copies the generated target-xlf name into column J
#   - stamps the "Latest Handback DateTime" in column K
#   - records the "handback file is not latest" warning in column P
#   - widens column P (Error Detail) to fit the warning text
$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/443d13d4ab8fa4837d2e56f0b3cff705396f036e/e2e/ed324dd0-a91d-4198-aeaf-d52b8e37e88d.md"
$handbackDisplay = "ed324dd0-a91d-4198-aeaf-d52b8e37e88d.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8d59920ce60ae1a63f80fd9d26773220774d57f/e2e/ed324dd0-a91d-4198-aeaf-d52b8e37e88d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/443d13d4ab8fa4837d2e56f0b3cff705396f036e/e2e/ed324dd0-a91d-4198-aeaf-d52b8e37e88d.md."

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns("P:P").ColumnWidth = 39.16
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsZh.Range("J8").Value = $wsZh.Range("G8").Text
$wsZh.Range("K8").Value = "2016-08-27 16:42:45"
$wsZh.Range("P8").Value = $errorDetail

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns("P:P").ColumnWidth = 39.16
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsDe.Range("J8").Value = $wsDe.Range("G8").Text
$wsDe.Range("K8").Value = "2016-08-27 16:42:51"
$wsDe.Range("P8").Value = $errorDetail
